# Rewrites the invoice data sheet:
#   * renames/reorders some headers and appends 5 new trailing columns
#     (price, taxAmount, taxRate, subTotal, items)
#   * replaces the 4 existing data rows with new invoice rows
#   * appends 2 brand-new invoice rows (rows 6 and 7)
#
# Text-look-alike values (dates, zero-padded numbers, phone numbers, "1, 1"
# style lists, "805.00" style amounts, ...) must be preserved verbatim as
# TEXT, not auto-converted to numbers/dates by Excel - so every text cell
# has its NumberFormat forced to "@" before the value is written. A handful
# of cells (grandPrice/duesAmount on rows 2-5) are genuine numbers and are
# written without the text format so they land as numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper: write one row starting at column A -----------------------
# $Cells is an ordered array of [ "ColumnLetter", value ] pairs. A bare
# number is written as-is so it lands as a real numeric cell. A string
# that *looks* numeric/date-like ("2025-08-16", "7894561230", "405.00", a
# lone "0", ...) would otherwise get silently reinterpreted by Excel as a
# number/date, so those get their cell NumberFormat forced to "@" (Text)
# before the value is written, which keeps the literal text intact.
# $null entries are skipped (cell stays/becomes empty).
# NOTE: named parameters (-Row / -Cells) are not reliably bound by this
# host's PowerShell interpreter, so this takes plain positional args.
function Set-RowCells {
    param($Row, $Cells)

    foreach ($pair in $Cells) {
        $col = $pair[0]
        $val = $pair[1]
        if ($null -eq $val) { continue }

        $ref = "$col$Row"
        if ($val -is [string] -and $val -match '^-?\d+(\.\d+)?$|^\d{4}-\d{2}-\d{2}$') {
            $ws.Range($ref).NumberFormat = "@"
        }
        $ws.Range($ref).Value = $val
    }
}

# --- header row ----------------------------------------------------------
Set-RowCells 1 @(
    ,@("A", "invoiceNumber")
    ,@("B", "date")
    ,@("C", "name")
    ,@("D", "mobileNumber")
    ,@("E", "product")
    ,@("F", "quantity")
    ,@("G", "originalPrice")
    ,@("H", "offerPrice")
    ,@("I", "grandPrice")
    ,@("J", "paidAmount")
    ,@("K", "duesAmount")
    ,@("L", "address")
    ,@("M", "bookingDate")
    ,@("N", "deliveryDate")
    ,@("O", "price")
    ,@("P", "taxAmount")
    ,@("Q", "taxRate")
    ,@("R", "subTotal")
    ,@("S", "items")
)

# --- row 2: INV-20250816-006 / Ram ---------------------------------------
Set-RowCells 2 @(
    ,@("A", "INV-20250816-006")
    ,@("B", "2025-08-16")
    ,@("C", "Ram")
    ,@("D", "7894561230")
    ,@("E", "zeiss dura vision blue protecgt, Sunglass")
    ,@("F", "1, 1")
    ,@("G", "3250, 700")
    ,@("H", "2800, 600")
    ,@("I", 3740)
    ,@("J", "2000")
    ,@("K", 1740)
    ,@("L", "Abc ")
    ,@("M", "2025-08-16")
    ,@("N", "2025-08-17")
)

# --- row 3: INV-20250816-005 / Amit Kumar --------------------------------
Set-RowCells 3 @(
    ,@("A", "INV-20250816-005")
    ,@("B", "2025-08-16")
    ,@("C", "Amit Kumar")
    ,@("D", "9876589870")
    ,@("E", "Sunglass, Red Frame")
    ,@("F", "2, 2")
    ,@("G", "700, 650")
    ,@("H", "600, 450")
    ,@("I", 2310)
    ,@("J", "2310")
    ,@("K", 0)
    ,@("L", "Eldeco Udhyaan - 02")
    ,@("M", "2025-08-16")
    ,@("N", "2025-08-18")
)

# --- row 4: INV-20250816-004 / Kunal -------------------------------------
Set-RowCells 4 @(
    ,@("A", "INV-20250816-004")
    ,@("B", "2025-08-16")
    ,@("C", "Kunal")
    ,@("D", "7894561230")
    ,@("E", "Lens, Frame")
    ,@("F", "1, 1")
    ,@("G", "900, 1000")
    ,@("H", "700, 700")
    ,@("I", 1540)
    ,@("J", "999")
    ,@("K", 541)
    ,@("L", "ldco")
    ,@("M", "2025-08-16")
    ,@("N", "2025-08-20")
)

# Old row 5 (the "Ajit / Frame, Lens" invoice) had G5/H5 populated; the new
# row 5 data below has no originalPrice/offerPrice list, so those two
# cells must be cleared out rather than overwritten.
$ws.Range("G5:H5").ClearContents()

# --- row 5: INV-20250816-003 / Keshav ------------------------------------
Set-RowCells 5 @(
    ,@("A", "INV-20250816-003")
    ,@("B", "2025-08-16")
    ,@("C", "Keshav")
    ,@("D", "7894561230")
    ,@("E", "sunglass")
    ,@("F", "1")
    ,@("I", 660)
    ,@("J", "500")
    ,@("K", 160)
    ,@("L", "ldco")
    ,@("M", "2025-08-16")
    ,@("N", "")
    ,@("O", "")
)

# --- row 6 (new): INV-20250816-002 / kUNAL -------------------------------
Set-RowCells 6 @(
    ,@("A", "INV-20250816-002")
    ,@("B", "2025-08-16")
    ,@("C", "kUNAL")
    ,@("D", "7894561230")
    ,@("E", "sunglass")
    ,@("F", "1")
    ,@("I", "805.00")
    ,@("J", "400")
    ,@("K", "405.00")
    ,@("L", "")
    ,@("M", "2025-08-16")
    ,@("N", "")
    ,@("O", "700")
    ,@("P", "105.00")
    ,@("Q", "15")
    ,@("R", "700.00")
)

# --- row 7 (new): INV-20250816-001 / Keshav ------------------------------
Set-RowCells 7 @(
    ,@("A", "INV-20250816-001")
    ,@("B", "2025-08-16")
    ,@("C", "Keshav")
    ,@("D", "7794561230")
    ,@("E", "sunglass")
    ,@("F", "1")
    ,@("I", "770.00")
    ,@("J", "500")
    ,@("K", "270.00")
    ,@("L", "")
    ,@("M", "2025-08-16")
    ,@("N", "")
    ,@("O", "700")
    ,@("P", "70.00")
    ,@("Q", "10")
    ,@("R", "700.00")
)
